$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.208510756492615
$ws.Range("B1").Value = 3.673158168792725
$ws.Range("C1").Value = 3.325173377990723
$ws.Range("D1").Value = 2.632070064544678
$ws.Range("E1").Value = 1.213170051574707
